$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns that carry row-specific data which gets swapped between rows.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

function Swap-RowValues($ws, $row1, $row2, $cols) {
    $vals1 = @{}
    $vals2 = @{}
    foreach ($col in $cols) {
        $vals1[$col] = $ws.Range("$col$row1").Value()
        $vals2[$col] = $ws.Range("$col$row2").Value()
    }
    foreach ($col in $cols) {
        $ws.Range("$col$row1").Value = $vals2[$col]
        $ws.Range("$col$row2").Value = $vals1[$col]
    }
}

# Row 35 <-> Row 40
Swap-RowValues $ws 35 40 $cols

# Row 38 <-> Row 39
Swap-RowValues $ws 38 39 $cols
